$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "- VIGEX_rsem_expression_results_genes_counts.tsv: ..."
# Merge the bookmark-separated trailing "." into the preceding sentence run
# so the paragraph reads "...external validation of VIGEX signature." as one
# run (no more separate "_GoBack" bookmark split in this paragraph).
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*VIGEX_rsem_expression_results_genes_counts.tsv*") {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $null = $targetPara.Range.Find.Execute(
        "of VIGEX signature" + [char]46,
        $true, $false, $false, $false, $false,
        $true, 1, $false,
        "of VIGEX signature.", 2)
}

# ---------------------------------------------------------------------------
# Change 2: insert a brand-new paragraph describing the TCGA files right
# after the "- additional_datasets: ..." paragraph (and before the trailing
# blank paragraphs). The "_GoBack" bookmark that used to sit inside the
# paragraph handled in Change 1 now lives inside this new paragraph instead.
# ---------------------------------------------------------------------------
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*additional_datasets:*") {
        $anchorPara = $p
    }
}

if ($anchorPara -ne $null) {
    $anchorPara.Range.InsertParagraphAfter()

    $newIndex = $anchorPara.Index + 1
    $newPara = $d.Paragraphs($newIndex)

    $tcgaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang/></w:rPr><w:t xml:space="preserve">- TCGA: contains expected counts (tcga_gene_expected_countlog2.coadmss.zip and tcga_gene_expected_countlog2.paad.zip for the Colon and Pancreas dataset respectively), the VIGEX score of the samples (COADMSS.tsv and PAAD.tsv </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default"/></w:rPr><w:t>for the Colon and Pancreas dataset respectively</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang/></w:rPr><w:t xml:space="preserve">) and the clinical data (COADREAD_MSS_survival.tsv and PAAD_survival.tsv </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default"/></w:rPr><w:t>for the Colon and Pancreas dataset respectively</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang/></w:rPr><w:t>).</w:t></w:r></w:p>
'@

    $null = $newPara.Range.InsertXML($tcgaXml)
}
